$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2025-2")

# Row 11: EP RIMAC
$ws.Range("A11").Value = "2025-2"
$ws.Range("B11").Value = "EP RIMAC"
$ws.Range("C11").Value = "Embarcación Pesquera"
$ws.Range("C2").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$ws.Range("D11").Value = "A.S/0045"
$ws.Range("E11").Value = "A.S/0045-225"

# Row 12
$ws.Range("A12").Value = "2025-2"
$ws.Range("C12").Value = "Embarcación Pesquera"
$ws.Range("C2").Copy()
$ws.Range("C12").PasteSpecial(-4122)
$ws.Range("D12").Value = "A.S/0046"
$ws.Range("E12").Value = "A.S/0046-225"
$ws.Range("B12").Value = "EP TIBURON 9"

$excel.CutCopyMode = $false
$ws.Range("C15").Select()
